# Generate Report for Archive
# - Status moves from "Ready for handoff" to "In Translation" on every sheet
#   that carries the per-file localization status (Overview!E2/F2, the
#   zh-cn sheet's Status cell, and the de-de sheet's Status cell).
# - The Status columns are narrower afterwards (re-fit to the new, shorter
#   text), so their column widths shrink to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text (was "Ready for handoff").
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns to fit the shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
